$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered save games) - row 2 values
$ws.Range("B2").Value = 0.1169995834814548
$ws.Range("C2").Value = 74547488392974520000000000.0
$ws.Range("D2").Value = 189.6080260415259
$ws.Range("E2").Value = 2521694498980204000000000000.0
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2596241987373178000000000000.0
